# CMP73010.docx edit: clean up the ">>> your stuff..." line, append a new
# paragraph with the "version control" blurb, and move the _GoBack
# bookmark into that new paragraph.

$d = $word.ActiveDocument

# --- 1. Clean up paragraph 4 ( ">>>  your stuff after this line >>>" ) ---
# In the original file this text is split across several runs with a
# <w:proofErr/> pair in the middle (Word's grammar-check markers). Rebuild
# it as a single plain run with no proofErr markers.
$p4 = $d.Paragraphs(4)
$p4Range = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$p4Range.Delete()
$p4Range.InsertBefore(">>>  your stuff after this line >>>")

# --- 2. Insert a brand new, empty paragraph right after paragraph 4 ---
$insPoint = $d.Paragraphs(4).Range
$insPoint.Collapse(0)
$insPoint.InsertParagraphAfter()

# The freshly created paragraph is now Paragraphs(5); "Ben changing
# things up!" has shifted down to Paragraphs(6). Since the new paragraph
# is empty, anchor the working range on its Start (its End is ambiguous
# with the start of the following paragraph).
$newPara = $d.Paragraphs(5)
$r = $d.Range($newPara.Range.Start, $newPara.Range.Start)

function Add-Chunk($range, $text, $bold) {
    $start = $range.End
    $range.InsertAfter($text)
    $chunk = $d.Range($start, $start + $text.Length)
    $chunk.Font.Name = "Arial"
    $chunk.Font.Color = 2236962
    if ($bold) {
        $chunk.Font.Bold = $true
    }
    return $chunk
}

Add-Chunk $r "A component of software configuration " $false | Out-Null
Add-Chunk $r "management" $true | Out-Null
Add-Chunk $r "," $false | Out-Null
Add-Chunk $r " version" $true | Out-Null
Add-Chunk $r " control" $true | Out-Null
Add-Chunk $r ", also known as revision " $false | Out-Null
Add-Chunk $r "control" $true | Out-Null
Add-Chunk $r " or " $false | Out-Null
Add-Chunk $r "source" $false | Out-Null
Add-Chunk $r " control" $true | Out-Null
Add-Chunk $r ", is the " $false | Out-Null
Add-Chunk $r "management" $true | Out-Null
Add-Chunk $r " of changes to documents, computer programs, large web sites, and other collectio" $false | Out-Null

# The _GoBack bookmark sits between "collectio" and "ns of information."
# in the final text, so mark that point before inserting the rest, then
# move Word's special _GoBack bookmark there (adding a bookmark named
# "_GoBack" relocates the existing one).
$bmPoint = $r.Duplicate
$bmPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

Add-Chunk $r "ns of information." $false | Out-Null
